$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data
$ws.Cells.Item(2,1).Value = '281474991153435-1752607648972'
$ws.Cells.Item(2,2).Value = 'No Seat Belt'
$ws.Cells.Item(2,3).Value = '2025-07-15T13:27:28.972'
$ws.Cells.Item(2,4).Value = "'281474991153435"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'141"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(2,6).Value = "'51834023"
$ws.Cells.Item(2,6).Style = "Normal"
$ws.Cells.Item(2,7).Value = 'DANIEL IÑIGUEZ'
$ws.Cells.Item(2,8).Value = 20.74040777
$ws.Cells.Item(2,9).Value = -103.385143889
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991153435/1752607646472/h9ZTBhLiGa-camera-video-segment-driver-1752607648972.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJ77774AX%2F20250716%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250716T140031Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC%2Bk16RMAd5LP46ZhpMPkPnCFOnCv2BiXmVyxE%2By6DRpAIhAM7SjCvOcLgaXX4Ia9yZ6PzVOFLcwVLIxlHwH2rguOrwKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgyWKLBLPu6B3ua7jDsqugPKA8vv8b5OdDs6paM8ljKlC4mB4Ln%2FnlzS%2B0%2FpTT9UfOpYgGfRKOMGFpeKFfGQVsl4iS8eOu3xOUZLwhHfVzUQPEK4css2Ef81qJNBvHFi%2B8XHIushSxbCtgdxGO4M2SISIPhN45OSXUdEiBi54NuyFUCLPUcKGhJZN8%2FMS2gI8Nzgg9VM1RqQEVJD2WmlDmbNZkzhdeUH0RS3kjxEimodDgrd8IT2xCIgAMv8QS51I7IxiIZI2TuEAaARk6YDLC2nTb%2BlQxfwrpk8HV49LoUnSZS0aCXEVcuC89Q1YhJyyF4pkBEaP2V%2FKnNw%2B4ddAO7MLwhaue0phSvzNZpAFrpmKupBesBJTZvVrZ6WY6eyBUwzUjB%2Fbbctq06enwnwwGA3fwnxEbFTx%2F8zDoGWtIvyBBKa7PVkQlnz6vCjAZZAkhVeGeRgB1WPdIs6H9qqwWt0iBjIhYhvceMawKsmdzZYmI3G2C24JIBlqfiTS%2FGWXmvFWgg49ovW69s%2BWkip9zll4hXOaz%2BvqWDHEn7jxlTKEX6Uree0lO%2FDxDUJhdfe4%2Bn5YrZ51%2FFzcscdg3r0P8odvKK0porOtvCAMMPD3sMGOqQBJnlWfp1UxCh6GxiE5PgWnGDdQtZg%2B7Kyzr0Wt00N3Ec3et6kmZXeBkU9zgEmo8jE5m2Be2YWTLUWutff8%2B13ZUMGwW7q5NNb0uUXXbmRYxZz3DtRvQHgpLFa4VwAQTaHEwT%2F7SYkaOUfopbHiu9W5NG0O%2F5q%2FByrE5H1JB0fHS6i0dO4LL27yMifSV4b3MnLE7TTc0iGcBR7bswnGP4vk93Dl6U%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2016%20Jul%202025%2022%3A00%3A31%20GMT&X-Amz-Signature=10efa38cf78937a85ce18127359873462d9af65cc01f3f8ac987ed26e1853b77'
$ws.Cells.Item(2,12).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991153435/1752607646472/aPa6fObqGp-camera-video-segment-1752607648972.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSJ77774AX%2F20250716%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250716T140031Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC%2Bk16RMAd5LP46ZhpMPkPnCFOnCv2BiXmVyxE%2By6DRpAIhAM7SjCvOcLgaXX4Ia9yZ6PzVOFLcwVLIxlHwH2rguOrwKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgyWKLBLPu6B3ua7jDsqugPKA8vv8b5OdDs6paM8ljKlC4mB4Ln%2FnlzS%2B0%2FpTT9UfOpYgGfRKOMGFpeKFfGQVsl4iS8eOu3xOUZLwhHfVzUQPEK4css2Ef81qJNBvHFi%2B8XHIushSxbCtgdxGO4M2SISIPhN45OSXUdEiBi54NuyFUCLPUcKGhJZN8%2FMS2gI8Nzgg9VM1RqQEVJD2WmlDmbNZkzhdeUH0RS3kjxEimodDgrd8IT2xCIgAMv8QS51I7IxiIZI2TuEAaARk6YDLC2nTb%2BlQxfwrpk8HV49LoUnSZS0aCXEVcuC89Q1YhJyyF4pkBEaP2V%2FKnNw%2B4ddAO7MLwhaue0phSvzNZpAFrpmKupBesBJTZvVrZ6WY6eyBUwzUjB%2Fbbctq06enwnwwGA3fwnxEbFTx%2F8zDoGWtIvyBBKa7PVkQlnz6vCjAZZAkhVeGeRgB1WPdIs6H9qqwWt0iBjIhYhvceMawKsmdzZYmI3G2C24JIBlqfiTS%2FGWXmvFWgg49ovW69s%2BWkip9zll4hXOaz%2BvqWDHEn7jxlTKEX6Uree0lO%2FDxDUJhdfe4%2Bn5YrZ51%2FFzcscdg3r0P8odvKK0porOtvCAMMPD3sMGOqQBJnlWfp1UxCh6GxiE5PgWnGDdQtZg%2B7Kyzr0Wt00N3Ec3et6kmZXeBkU9zgEmo8jE5m2Be2YWTLUWutff8%2B13ZUMGwW7q5NNb0uUXXbmRYxZz3DtRvQHgpLFa4VwAQTaHEwT%2F7SYkaOUfopbHiu9W5NG0O%2F5q%2FByrE5H1JB0fHS6i0dO4LL27yMifSV4b3MnLE7TTc0iGcBR7bswnGP4vk93Dl6U%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2016%20Jul%202025%2022%3A00%3A31%20GMT&X-Amz-Signature=3f1232eebec9b0cbdec98e3b32884fa68172045606db7cc78c8bd0e9c3acc6e7'

# Row 3 data
$ws.Cells.Item(3,1).Value = '281474991395157-1752595308013'
$ws.Cells.Item(3,2).Value = 'Harsh Brake'
$ws.Cells.Item(3,3).Value = '2025-07-15T10:01:48.013'
$ws.Cells.Item(3,4).Value = "'281474991395157"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'126"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,6).Value = "'53418410"
$ws.Cells.Item(3,6).Style = "Normal"
$ws.Cells.Item(3,7).Value = 'FELIPE'
$ws.Cells.Item(3,8).Value = 20.68313656
$ws.Cells.Item(3,9).Value = -103.36120946
$ws.Cells.Item(3,10).Value = 0.8949801921844482
$ws.Cells.Item(3,11).Value = 'No video URL'
$ws.Cells.Item(3,12).Value = 'No video URL'

# Row 4 data
$ws.Cells.Item(4,1).Value = '281474991395157-1752593236230'
$ws.Cells.Item(4,2).Value = 'Harsh Brake'
$ws.Cells.Item(4,3).Value = '2025-07-15T09:27:16.230'
$ws.Cells.Item(4,4).Value = "'281474991395157"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'126"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(4,6).Value = "'53418410"
$ws.Cells.Item(4,6).Style = "Normal"
$ws.Cells.Item(4,7).Value = 'FELIPE'
$ws.Cells.Item(4,8).Value = 20.65680467
$ws.Cells.Item(4,9).Value = -103.34967407
$ws.Cells.Item(4,10).Value = 0.8669819235801697
$ws.Cells.Item(4,11).Value = 'No video URL'
$ws.Cells.Item(4,12).Value = 'No video URL'

